# Update the "想去人数" (interest count) figures that changed between data
# refreshes for the events that appear on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F9").Value  = 81
$wsExpo.Range("F10").Value = 290
$wsExpo.Range("F13").Value = 12037
$wsExpo.Range("F14").Value = 12480

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 81
$wsAll.Range("F11").Value = 290
$wsAll.Range("F14").Value = 12037
$wsAll.Range("F15").Value = 12480
